$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24; existing rows 24-45 shift down to 25-46.
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with this week's record (same market/product metadata
# as the surrounding rows, new date + price figures).
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C24").Value = "Arica y Parinacota"
$ws.Range("D24").Value = 45096
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = 100112003
$ws.Range("G24").Value = "Ajo"
$ws.Range("H24").Value = "Chino"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 600
$ws.Range("K24").Value = 17000
$ws.Range("L24").Value = 18000
$ws.Range("M24").Value = 17417
$ws.Range("N24").Value = "$/caja 10 kilos"
$ws.Range("O24").Value = "China"
$ws.Range("P24").Value = 1742
$ws.Range("Q24").Value = 10
$ws.Range("R24").Value = "Hortaliza"
